$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "moved J11 to the same 6pin connector" -- J11 is now populated on the same
# 6-position header (J7/J8/J10/J13/J14) as the other connectors, and its
# mating pin header (P7/P8/P10/P13/P14) likewise. That grows both RefDes
# lists to include J11 / P11, and the populated-position count for that
# connector goes from 5 to 6 (which flows into the two downstream rows that
# derive their quantities from it).
$ws.Range("E16").Value = "J7 J8 J10 J11 J13 J14"
$ws.Range("E17").Value = "P7 P8 P10 P11 P13 P14"
$ws.Range("F16").Value = 6

# Reflect where the editor ended up looking at this change: zoomed in to
# 130% and scrolled/selected near the edited rows.
$win = $excel.ActiveWindow
$win.Zoom = 130
$ws.Range("E19").Select()
